$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.166.70'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '2.376.58'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'549.52"
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = "'138.78"
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('D9').Value = '2.377.63'
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('E10').Value = '  +2.52%  '
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').Value = "'25.04"
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').Value = '2.787.97'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').Value = '61.093.83'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').Value = '2.389.21'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = "'10.84"
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'321.90"
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = "'4.16"
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').Value = "'6.73"
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('D24').Value = "'64.36"
$ws.Range('E24').Value = '  +1.33%  '
$ws.Range('D25').Value = "'1.70"
$ws.Range('E25').Value = '  -11.07%  '
$ws.Range('D26').Value = "'8.37"
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('D27').Value = "'1.00"
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '2.487.21'
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('D29').Value = "'8.17"
$ws.Range('E29').Value = '  +1.09%  '
$ws.Range('D30').Value = "'508.67"
$ws.Range('E30').Value = '  -3.15%  '
$ws.Range('E31').Value = '  -4.35%  '
$ws.Range('E32').Value = '  +3.10%  '
$ws.Range('E33').Value = '  -3.45%  '
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('E35').Value = '  -3.55%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = "'5.39"
$ws.Range('E39').Value = '  -2.13%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = "'1.87"
$ws.Range('E40').Value = '  +2.14%  '
$ws.Range('D41').Value = "'18.59"
$ws.Range('E41').Value = '  +2.78%  '
$ws.Range('D42').Value = "'146.66"
$ws.Range('E42').Value = '  +5.97%  '
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').Value = "'41.33"
$ws.Range('E44').Value = '  +2.55%  '
$ws.Range('D45').Value = "'148.40"
$ws.Range('E45').Value = '  +5.47%  '
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('E47').Value = '  -2.77%  '
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('D49').Value = "'19.31"
$ws.Range('E49').Value = '  -4.03%  '
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').Value = "'0.0910"
$ws.Range('E51').Value = '  +0.56%  '
